$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Section" and "Term" columns (old C:D)
$ws.Columns("C:D").Delete()

# Remove the extra data rows, keep only the header + first data row
$ws.Rows("3:5").Delete()

# Update the remaining data row's values (Days / Start Time / End Time)
$ws.Range("D2").Value = "M"
$ws.Range("E2").Value = 0.35416666666666669
$ws.Range("F2").Value = 0.4375

# Re-format the Start Time / End Time headers so they no longer look like
# the bold/wrapped header style
$ws.Range("A2").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Faculty name data cell drops the time number format, keeping the plain
# centred border style
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# Start / End time data cells keep (re-apply) the time number format
$ws.Range("E2:F2").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4142

$excel.CutCopyMode = 0

# Autofit the Start Time / End Time columns to their new (narrower) content
$ws.Columns("E:F").AutoFit()

# Shrink the data validation range down to the cells that remain
$ws.Range("B2:C5,G2:I5").Validation.Delete()
$ws.Range("B2,G2").Validation.Add(3, 1, 1, "#REF!")

# Update the selected cell shown when the workbook is opened
$ws.Range("F10").Select()

# Switch the page to portrait orientation
$ws.PageSetup.Orientation = 1
